$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update HTTPS Request joule measurements (Nexus 5X) in column B
$ws.Range("B2").Value = 135.363888
$ws.Range("B3").Value = 121.996511999999
$ws.Range("B4").Value = 120.76732800000001
$ws.Range("B5").Value = 119.23084799999999
$ws.Range("B6").Value = 119.384496
$ws.Range("B7").Value = 119.384496
$ws.Range("B8").Value = 119.384496
$ws.Range("B9").Value = 119.999088
$ws.Range("B10").Value = 119.384496
$ws.Range("B11").Value = 119.53814399999899
$ws.Range("B12").Value = 119.69179200000001
$ws.Range("B13").Value = 117.54071999999999
$ws.Range("B14").Value = 137.05401599999999
$ws.Range("B15").Value = 136.74672000000001
$ws.Range("B16").Value = 136.74672000000001
$ws.Range("B17").Value = 120.152736
$ws.Range("B18").Value = 119.23084799999999
$ws.Range("B19").Value = 119.384496
$ws.Range("B20").Value = 120.460032
$ws.Range("B21").Value = 119.384496
$ws.Range("B22").Value = 119.84544
$ws.Range("B23").Value = 119.0772
$ws.Range("B24").Value = 107.860895999999
$ws.Range("B25").Value = 119.53814399999899
$ws.Range("B26").Value = 134.90294399999999
$ws.Range("B27").Value = 135.67118400000001
$ws.Range("B28").Value = 134.288352
$ws.Range("B29").Value = 134.59564800000001
$ws.Range("B30").Value = 118.769904
$ws.Range("B31").Value = 119.69179200000001
$ws.Range("B32").Value = 116.223912
$ws.Range("B33").Value = 122.80262399999999
$ws.Range("B34").Value = 119.983176
$ws.Range("B35").Value = 120.60972
$ws.Range("B36").Value = 120.92299199999999
$ws.Range("B37").Value = 120.139811999999
$ws.Range("B38").Value = 120.296448
$ws.Range("B39").Value = 121.3929
$ws.Range("B40").Value = 119.669904
$ws.Range("B41").Value = 119.82653999999999
$ws.Range("B42").Value = 120.139811999999
$ws.Range("B43").Value = 121.549536
$ws.Range("B44").Value = 135.64677599999999
$ws.Range("B45").Value = 139.249404
$ws.Range("B46").Value = 135.49014
$ws.Range("B47").Value = 135.17686800000001
$ws.Range("B48").Value = 137.36977199999899
$ws.Range("B49").Value = 136.116683999999
$ws.Range("B50").Value = 135.33350399999901
$ws.Range("B51").Value = 134.70695999999899
$ws.Range("B52").Value = 135.80341200000001
$ws.Range("B53").Value = 135.64677599999999
$ws.Range("B54").Value = 129.06806399999999
$ws.Range("B55").Value = 135.960048
$ws.Range("B56").Value = 136.27331999999899
$ws.Range("B57").Value = 134.393688
$ws.Range("B58").Value = 119.82653999999999
$ws.Range("B59").Value = 119.199996
$ws.Range("B60").Value = 120.296448
$ws.Range("B61").Value = 120.453084
$ws.Range("B62").Value = 166.80923999999999
$ws.Range("B63").Value = 120.74540399999999
$ws.Range("B64").Value = 118.75593600000001
$ws.Range("B65").Value = 119.674151999999
$ws.Range("B66").Value = 118.44986400000001
$ws.Range("B67").Value = 119.521115999999
$ws.Range("B68").Value = 119.674151999999
$ws.Range("B69").Value = 118.143792
$ws.Range("B70").Value = 119.215043999999
$ws.Range("B71").Value = 118.60290000000001
$ws.Range("B72").Value = 123.500052
$ws.Range("B73").Value = 117.83772
$ws.Range("B74").Value = 136.355076
$ws.Range("B75").Value = 135.43686
$ws.Range("B76").Value = 136.20204000000001
$ws.Range("B77").Value = 135.89596800000001
$ws.Range("B78").Value = 137.579364
$ws.Range("B79").Value = 136.81418399999899
$ws.Range("B80").Value = 135.130788
$ws.Range("B81").Value = 136.049004
$ws.Range("B82").Value = 135.742932
$ws.Range("B83").Value = 136.355076
$ws.Range("B84").Value = 127.478988
$ws.Range("B85").Value = 136.661148
$ws.Range("B86").Value = 136.355076
$ws.Range("B87").Value = 135.89596800000001
$ws.Range("B88").Value = 135.130788
$ws.Range("B89").Value = 136.049004
$ws.Range("B90").Value = 135.43686
$ws.Range("B91").Value = 135.89596800000001

# Restore the last selection/scroll position recorded when the file was saved
$ws.Range("G84").Select()
